$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 2 new rows right after the last data row (row 27), before the blank spacer rows ---
$ws.Rows("28:29").Insert()

# The newly-inserted row 29 becomes the new "bottom of table" row, so it should carry the
# special bottom-border formatting that row 27 currently has (copy it before row 27 is changed).
$ws.Range("B27:J27").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 27 is no longer the last row of the table, so it now takes the regular "middle" row style
# (matching row 26, the row above it).
$ws.Range("B26:J26").Copy()
$ws.Range("B27:J27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The other new row (28) is also a regular "middle" row.
$ws.Range("B26:J26").Copy()
$ws.Range("B28:J28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Fill in the two new employees' data ---
$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1027950315"
$ws.Range("D28").Value = "REVELE BANGUERA PESTANA"
$ws.Range("E28").Value = "2509"
$ws.Range("F28").Value = 56940
$ws.Range("G28").Value = 1423500

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1048454047"
$ws.Range("D29").Value = "ENIS MERCEDEZ BLANCO ZUÑIGA"
$ws.Range("E29").Value = "2509"
$ws.Range("F29").Value = 56940
$ws.Range("G29").Value = 1423500

# --- 3. Update the summary figures at the top of the report ---
$ws.Range("E11").Value = 618685
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 13
